$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-14 (B:F) with the new Air Conditioning System test cases ---
$ws.Range("B2").Value = 'Default Temprature '
$ws.Range("C2").Value = 'Intial Temprature = 20'
$ws.Range("D2").Value = 'Intial Temprature = 20'
$ws.Range("E2").Value = 'Start With Temprature = 20'
$ws.Range("F2").Value = 'Pass'
$ws.Range("B3").Value = 'Incremant Button In Set Mode'
$ws.Range("C3").Value = 'Push The Incremant Button'
$ws.Range("D3").Value = 'Incremant In Temprature'
$ws.Range("E3").Value = 'Temprature Increment'
$ws.Range("F3").Value = 'Pass'
$ws.Range("B4").Value = 'Decremant Button In Set Mode'
$ws.Range("C4").Value = 'Push The Decremant Button'
$ws.Range("D4").Value = ' Decremant In Temprature'
$ws.Range("E4").Value = 'Temprature Decrement'
$ws.Range("F4").Value = 'pass'
$ws.Range("B5").Value = 'Set Button In Set Mode'
$ws.Range("C5").Value = 'Push The Set Button'
$ws.Range("D5").Value = 'Set New Temprature'
$ws.Range("E5").Value = 'User Can Set New Temprature'
$ws.Range("F5").Value = 'Pass'
$ws.Range("B6").Value = 'Adjust Button In Set Mode'
$ws.Range("C6").Value = 'Push The Adjust Button'
$ws.Range("D6").Value = 'Nothing'
$ws.Range("E6").Value = 'Nothing'
$ws.Range("F6").Value = 'pass'
$ws.Range("B7").Value = 'Reset Button In Set Mode'
$ws.Range("C7").Value = 'Push The Reset Button'
$ws.Range("D7").Value = 'Nothing'
$ws.Range("E7").Value = 'Nothing'
$ws.Range("F7").Value = 'pass'
$ws.Range("B8").Value = 'Buzzer Turn On'
$ws.Range("C8").Value = 'Temprature Exceed The Set Degree'
$ws.Range("D8").Value = 'Buzzer Turn On'
$ws.Range("E8").Value = 'Buzzer Turn On'
$ws.Range("F8").Value = 'pass'
$ws.Range("B9").Value = 'Buzzer Turn Off'
$ws.Range("C9").Value = 'Temprature Not Exceed The Set Degree'
$ws.Range("D9").Value = 'Buzzer turn Off'
$ws.Range("E9").Value = 'Buzzer Turn Off'
$ws.Range("F9").Value = 'pass'
$ws.Range("B10").Value = 'Buzzer Draw'
$ws.Range("C10").Value = 'Buzzer Char Draw On The Lcd When Buzzer Turn On'
$ws.Range("D10").Value = 'Buzzer Char Draw On Lcd'
$ws.Range("E10").Value = 'Buzzer Drawed'
$ws.Range("F10").Value = 'pass'
$ws.Range("B11").Value = 'Buzzer Clear'
$ws.Range("C11").Value = 'Buzzer Char Clear From The Lcd When Buzzer Turn OFF'
$ws.Range("D11").Value = 'Buzzer Char Clear From Lcd'
$ws.Range("E11").Value = 'Buzzer Cleared'
$ws.Range("F11").Value = 'pass'
$ws.Range("B12").Value = 'Current Temprature '
$ws.Range("C12").Value = 'Current Temprature Display On Lcd'
$ws.Range("D12").Value = 'Show Current Temprature On Lcd'
$ws.Range("E12").Value = 'Current Temprature Showed On Lcd'
$ws.Range("F12").Value = 'pass'
$ws.Range("B13").Value = 'Incremant Button In Idle Mode'
$ws.Range("C13").Value = 'Push Incremant Button'
$ws.Range("D13").Value = 'Error Massage On Lcd'
$ws.Range("E13").Value = 'Error Massage On Lcd'
$ws.Range("F13").Value = 'pass'
$ws.Range("B14").Value = 'Decremant Button In Idle Mode'
$ws.Range("C14").Value = 'Push Decremant Button'
$ws.Range("D14").Value = 'Error Massage On Lcd'
$ws.Range("E14").Value = 'Error Massage On Lcd'
$ws.Range("F14").Value = 'pass'
$ws.Range("A15").Value = 15
$ws.Range("B15").Value = 'Set Button In Idle Mode'
$ws.Range("C15").Value = 'Push Set Button'
$ws.Range("D15").Value = 'Error Massage On Lcd'
$ws.Range("E15").Value = 'Error Massage On Lcd'
$ws.Range("F15").Value = 'pass'
$ws.Range("A16").Value = 16
$ws.Range("B16").Value = 'Reset Button In Idle Mode'
$ws.Range("C16").Value = 'Push Reset Button'
$ws.Range("D16").Value = 'Temprature Return To The Default = 20'
$ws.Range("E16").Value = 'Temprature Return To The Default = 20'
$ws.Range("F16").Value = 'pass'
$ws.Range("A17").Value = 17
$ws.Range("B17").Value = 'Adjust Button In Idle Mode'
$ws.Range("C17").Value = 'Push Adjust Button'
$ws.Range("D17").Value = 'Turn To Set Mode'
$ws.Range("E17").Value = 'Turn To Set Mode'
$ws.Range("F17").Value = 'pass'
$ws.Range("A18").Value = 19
$ws.Range("B18").Value = 'Time Out'
$ws.Range("C18").Value = 'User Do Not Select Any Thing'
$ws.Range("D18").Value = 'Turn To Idle Mode With Default Temprature'
$ws.Range("E18").Value = 'Turn To Idle Mode With Default Temprature'
$ws.Range("F18").Value = 'Pass'

# --- Column widths (closest achievable to the target OOXML widths through the
#     ColumnWidth COM property, which this host rounds to 1/6-character steps) ---
$ws.Columns.Item(2).ColumnWidth = 41
$ws.Columns.Item(4).ColumnWidth = 37.3
$ws.Columns.Item(5).ColumnWidth = 36.5

# --- Restore the active selection left by the author on D23 ---
$ws.Range("D23").Select()
